# Last of the localization and gui fixes
#
# Applies:
#  - "Merge" sheet (sheet1): new "P"/"NO NEED" marker cells in column B (and
#    one C cell), re-using existing fill styles; view selection/scroll update.
#  - "Problems" sheet (sheet2): new "NO NEED" marker + new "UNSURE IF ISSUE"
#    note; view selection update.
#  - "Sheet5" (sheet5): three new blank, explicitly no-fill cells in column B;
#    view selection update.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122
$xlNone = -4142

$wsMerge = $wb.Worksheets.Item("Merge")
$wsProblems = $wb.Worksheets.Item("Problems")
$wsSheet5 = $wb.Worksheets.Item("Sheet5")

# ---------------------------------------------------------------------------
# Merge sheet: stamp column B (and C50) with the "P" (shared string 11) /
# "NO NEED" markers. Re-use the formatting already present on existing
# marker cells so the workbook keeps reusing the same style + shared string
# records instead of minting new ones.
# ---------------------------------------------------------------------------

# Reference cells already carrying the two fills used below.
$pStyleRef = $wsMerge.Range("B2")         # fillId 2 ("P" green)
$noNeedStyleRef = $wsMerge.Range("B68")   # fillId 2 ("NO NEED")

$pTargets = @("B19", "B21", "B23", "B24", "B37", "C50", "B53", "B66")
foreach ($addr in $pTargets) {
    $pStyleRef.Copy()
    $wsMerge.Range($addr).PasteSpecial($xlPasteFormats)
    $wsMerge.Range($addr).Value = "P"
}

$noNeedTargets = @("B28", "B35", "B36", "B43", "B44", "B45", "B46", "B47", "B51", "B54")
foreach ($addr in $noNeedTargets) {
    $noNeedStyleRef.Copy()
    $wsMerge.Range($addr).PasteSpecial($xlPasteFormats)
    $wsMerge.Range($addr).Value = "NO NEED"
}

$excel.CutCopyMode = $false

# View: scroll position + active selection on the Merge sheet.
$wsMerge.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsMerge.Range("D20").Select()

# ---------------------------------------------------------------------------
# Problems sheet: one more "NO NEED" marker, plus a brand new note cell.
# ---------------------------------------------------------------------------

$noNeedStyleRef.Copy()
$wsProblems.Range("B4").PasteSpecial($xlPasteFormats)
$wsProblems.Range("B4").Value = "NO NEED"

$maybeStyleRef = $wsMerge.Range("B48")   # fillId 5 ("MAYBE"-style note fill)
$maybeStyleRef.Copy()
$wsProblems.Range("B5").PasteSpecial($xlPasteFormats)
$wsProblems.Range("B5").Value = "UNSURE IF ISSUE"

$excel.CutCopyMode = $false

$wsProblems.Activate()
$wsProblems.Range("D9").Select()

# ---------------------------------------------------------------------------
# Sheet5: three explicitly "no fill" blank cells in column B.
# ---------------------------------------------------------------------------

foreach ($addr in @("B4", "B5", "B6")) {
    $cell = $wsSheet5.Range($addr)
    $cell.Interior.ColorIndex = 6
    $cell.Interior.ColorIndex = $xlNone
}

$wsSheet5.Activate()
$wsSheet5.Range("D9").Select()
